# Convert the "I" column from milliamps to amps:
#  - rename the header in B1 from "I [mA]" to "I [A]"
#  - rescale every data value in column B (rows 2-20) by dividing by 1000

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the column header
$ws.Range("B1").Value = "I [A]"

# New (already-converted) current values, in Amps, for rows 2 through 20
$newValues = @{
    2  = 0.0002512480738027979
    3  = 0.000502021016818301
    4  = 0.0007518481776414415
    5  = 0.001000267765505563
    6  = 0.001246831039134866
    7  = 0.001491106216376509
    8  = 0.001732682028419033
    9  = 0.001971170855562723
    10 = 0.002206211396452011
    11 = 0.002437470838597086
    12 = 0.004491102074175144
    13 = 0.006020330675915344
    14 = 0.007089765698382925
    15 = 0.007824789858269987
    16 = 0.008334021839423876
    17 = 0.00869371147732024
    18 = 0.008953709900542435
    19 = 0.009146065885575001
    20 = 0.009291520335781388
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}
